$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new column width for column CP (94) to match existing data columns (width 12)
$ws.Range("CP1").ColumnWidth = 11.17

# Row 1: header date, stored as text (like the other date header cells in row 1)
$ws.Range("CP1").NumberFormat = "@"
$ws.Range("CP1").Value2 = "2024/12/11"
$ws.Range("CO1").Copy()
$ws.Range("CP1").PasteSpecial(-4122)

# Data rows 2-53: set values
$ws.Range("CP2").Value2 = 195.6
$ws.Range("CP3").Value2 = 119.7
$ws.Range("CP4").Value2 = 140
$ws.Range("CP5").Value2 = 249
$ws.Range("CP6").Value2 = 162.9
$ws.Range("CP7").Value2 = 216.1
$ws.Range("CP8").Value2 = 170.2
$ws.Range("CP9").Value2 = 123.8
$ws.Range("CP10").Value2 = 145.3
$ws.Range("CP11").Value2 = 175.8
$ws.Range("CP12").Value2 = 138.6
$ws.Range("CP13").Value2 = 145.9
$ws.Range("CP14").Value2 = 203.5
$ws.Range("CP15").Value2 = 147.7
$ws.Range("CP16").Value2 = 137.4
$ws.Range("CP17").Value2 = 128.8
$ws.Range("CP18").Value2 = 164.6
$ws.Range("CP19").Value2 = 249.8
$ws.Range("CP20").Value2 = 130.6
$ws.Range("CP21").Value2 = 116.8
$ws.Range("CP22").Value2 = 131
$ws.Range("CP23").Value2 = 109.9
$ws.Range("CP24").Value2 = 206.1
$ws.Range("CP25").Value2 = 130.8
$ws.Range("CP26").Value2 = 126.9
$ws.Range("CP27").Value2 = 191.9
$ws.Range("CP28").Value2 = 132.6
$ws.Range("CP29").Value2 = 124
$ws.Range("CP30").Value2 = 159.8
$ws.Range("CP31").Value2 = 157.5
$ws.Range("CP32").Value2 = 141.2
$ws.Range("CP33").Value2 = 129.9
$ws.Range("CP34").Value2 = 422.3
$ws.Range("CP35").Value2 = 171.4
$ws.Range("CP36").Value2 = 127.8
$ws.Range("CP37").Value2 = 141.2
$ws.Range("CP38").Value2 = 174.1
$ws.Range("CP39").Value2 = 125.7
$ws.Range("CP40").Value2 = 149.8
$ws.Range("CP41").Value2 = 129.8
$ws.Range("CP42").Value2 = 122.1
$ws.Range("CP43").Value2 = 171.8
$ws.Range("CP44").Value2 = 175.9
$ws.Range("CP45").Value2 = 162.5
$ws.Range("CP46").Value2 = 302.8
$ws.Range("CP47").Value2 = 136.6
$ws.Range("CP48").Value2 = 168.4
$ws.Range("CP49").Value2 = 141.1
$ws.Range("CP50").Value2 = 129.3
$ws.Range("CP51").Value2 = 213.5
$ws.Range("CP52").Value2 = 165.3
$ws.Range("CP53").Value2 = 127.5

# Apply matching cell styles (reusing existing style indices) per group
$ws.Range("A2").Copy()
foreach ($cell in $ws.Range("CP2,CP4,CP5,CP6,CP7,CP8,CP10,CP11,CP13,CP14,CP15,CP18,CP19,CP24,CP27,CP30,CP31,CP32,CP34,CP35,CP37,CP38,CP40,CP43,CP44,CP45,CP46,CP48,CP49,CP51,CP52").Cells) {
    $cell.PasteSpecial(-4122)
}
$ws.Range("D2").Copy()
foreach ($cell in $ws.Range("CP3,CP9,CP21,CP23,CP29,CP42").Cells) {
    $cell.PasteSpecial(-4122)
}
$ws.Range("N2").Copy()
foreach ($cell in $ws.Range("CP12,CP16,CP17,CP20,CP22,CP25,CP26,CP28,CP33,CP36,CP39,CP41,CP47,CP50,CP53").Cells) {
    $cell.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
